$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "theta_threshold_range" row (row 5); this shifts the
# "pie_threshold_range" row (formerly row 6) up into row 5, and the
# now-unreferenced "theta_threshold_range" shared string is dropped
# automatically on save.
$ws.Rows("5").Delete()

# The cell that used to be B6 (style carried the Times New Roman font,
# s="3") is now B5; restore it to the plain numeric style used by the
# rest of column B/C (same style as B2, s="2").
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false | Out-Null

# Update the numeric values per the new measurements.
$ws.Range("B2").Value = 4.5
$ws.Range("C2").Value = 11.4

$ws.Range("B3").Value = 4.7
$ws.Range("C3").Value = 10.4

$ws.Range("B4").Value = 0.8
$ws.Range("C4").Value = 1.4

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Update the selected cell shown when the sheet is reopened.
$ws.Range("C3").Select() | Out-Null

# Page setup now prints on A4 portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1 | Out-Null
